$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-27 06:18:32'
$ws.Range("E3").Value = '2026-02-27 06:18:34'
$ws.Range("O3").Value = '3.0 °C'
$ws.Range("E4").Value = '2026-02-27 06:18:37'
$ws.Range("O4").Value = '6.8 °C'
$ws.Range("E5").Value = '2026-02-27 06:18:40'
$ws.Range("E6").Value = '2026-02-27 06:18:42'
$ws.Range("I6").Value = '0.1 mm'
$ws.Range("L6").Value = '6.1 km/h - 29º 5:38 TU'
$ws.Range("E7").Value = '2026-02-27 06:18:44'
$ws.Range("J7").Value = '1025.9 hPa'
$ws.Range("K7").Value = '-0.1 MJ/m2'
$ws.Range("N7").Value = '8.9 °C 5:57 TU'
$ws.Range("O7").Value = '10.0 °C'
$ws.Range("E8").Value = '2026-02-27 06:18:47'
$ws.Range("J8").Value = '1025.3 hPa'
$ws.Range("O8").Value = '12.0 °C'
$ws.Range("E9").Value = '2026-02-27 06:18:49'
$ws.Range("M9").Value = '9.2 °C 5:58 TU'
$ws.Range("O9").Value = '8.5 °C'
$ws.Range("E10").Value = '2026-02-27 06:18:52'
$ws.Range("N10").Value = '7.6 °C 5:59 TU'
$ws.Range("O10").Value = '9.1 °C'
$ws.Range("E11").Value = '2026-02-27 06:18:54'
$ws.Range("N11").Value = '1.1 °C 5:47 TU'
$ws.Range("O11").Value = '2.3 °C'
$ws.Range("E12").Value = '2026-02-27 06:18:57'
$ws.Range("M12").Value = '9.9 °C 5:57 TU'
$ws.Range("O12").Value = '8.3 °C'
$ws.Range("E13").Value = '2026-02-27 06:18:58'
$ws.Range("N13").Value = '-3.4 °C 5:59 TU'
$ws.Range("O13").Value = '-1.5 °C'
$ws.Range("E14").Value = '2026-02-27 06:18:59'
$ws.Range("N14").Value = '6.2 °C 5:36 TU'
$ws.Range("O14").Value = '8.6 °C'
$ws.Range("E15").Value = '2026-02-27 06:19:00'
$ws.Range("M15").Value = '9.7 °C 5:53 TU'
$ws.Range("O15").Value = '8.5 °C'
$ws.Range("E16").Value = '2026-02-27 06:19:01'
$ws.Range("E17").Value = '2026-02-27 06:19:03'
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = '33%'
$ws.Range("L17").Value = '45.7 km/h - 229º 5:12 TU'
$ws.Range("N17").Value = '6.2 °C 5:44 TU'
$ws.Range("O17").Value = '7.1 °C'
$ws.Range("E18").Value = '2026-02-27 06:19:04'
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = '96%'
$ws.Range("J18").Value = '1025.8 hPa'
$ws.Range("O18").Value = '9.7 °C'
$ws.Range("E19").Value = '2026-02-27 06:19:05'
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = '83%'
$ws.Range("O19").Value = '7.8 °C'
$ws.Range("E20").Value = '2026-02-27 06:19:06'
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = '59%'
$ws.Range("O20").Value = '2.0 °C'
$ws.Range("E21").Value = '2026-02-27 06:19:07'
$ws.Range("N21").Value = '1.8 °C 5:47 TU'
$ws.Range("O21").Value = '3.6 °C'
$ws.Range("E22").Value = '2026-02-27 06:19:08'
$ws.Range("E23").Value = '2026-02-27 06:19:11'
$ws.Range("O23").Value = '2.7 °C'
$ws.Range("E24").Value = '2026-02-27 06:19:13'
$ws.Range("O24").Value = '4.1 °C'
$ws.Range("E25").Value = '2026-02-27 06:19:16'
$ws.Range("E26").Value = '2026-02-27 06:19:18'
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = '44%'
$ws.Range("J26").Value = '1024.6 hPa'
$ws.Range("E27").Value = '2026-02-27 06:19:21'
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = '45%'
$ws.Range("E28").Value = '2026-02-27 06:19:23'
$ws.Range("O28").Value = '5.5 °C'
$ws.Range("E29").Value = '2026-02-27 06:19:26'
$ws.Range("E30").Value = '2026-02-27 06:19:28'
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = '99%'
$ws.Range("N30").Value = '9.0 °C 5:33 TU'
$ws.Range("O30").Value = '9.8 °C'
$ws.Range("E31").Value = '2026-02-27 06:19:30'
$ws.Range("E32").Value = '2026-02-27 06:19:33'
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = '92%'
$ws.Range("N32").Value = '-1.2 °C 5:47 TU'
$ws.Range("O32").Value = '0.9 °C'
$ws.Range("E33").Value = '2026-02-27 06:19:35'
$ws.Range("N33").Value = '0.3 °C 5:56 TU'
$ws.Range("O33").Value = '2.3 °C'
$ws.Range("E34").Value = '2026-02-27 06:19:37'
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = '45%'
$ws.Range("M34").Value = '4.4 °C 5:57 TU'
$ws.Range("O34").Value = '2.2 °C'
$ws.Range("E35").Value = '2026-02-27 06:19:40'
$ws.Range("N35").Value = '8.5 °C 5:59 TU'
$ws.Range("O35").Value = '10.0 °C'
$ws.Range("E36").Value = '2026-02-27 06:19:42'
$ws.Range("J36").Value = '1025.9 hPa'
$ws.Range("M36").Value = '10.4 °C 5:45 TU'
$ws.Range("O36").Value = '9.3 °C'
$ws.Range("E37").Value = '2026-02-27 06:19:45'
$ws.Range("J37").Value = '1028.5 hPa'
$ws.Range("N37").Value = '1.6 °C 5:59 TU'
$ws.Range("E38").Value = '2026-02-27 06:19:47'
$ws.Range("N38").Value = '6.6 °C 5:39 TU'
$ws.Range("O38").Value = '7.5 °C'
$ws.Range("E39").Value = '2026-02-27 06:19:50'
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = '20%'
$ws.Range("O39").Value = '5.1 °C'
$ws.Range("E40").Value = '2026-02-27 06:19:52'
$ws.Range("N40").Value = '0.7 °C 5:57 TU'
$ws.Range("O40").Value = '1.9 °C'
$ws.Range("E41").Value = '2026-02-27 06:19:55'
$ws.Range("J41").Value = '1025.8 hPa'
$ws.Range("N41").Value = '6.6 °C 5:59 TU'
$ws.Range("O41").Value = '8.6 °C'
$ws.Range("E42").Value = '2026-02-27 06:19:57'
$ws.Range("E43").Value = '2026-02-27 06:20:00'
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = '100%'
$ws.Range("N43").Value = '1.7 °C 5:57 TU'
$ws.Range("O43").Value = '3.9 °C'
$ws.Range("E44").Value = '2026-02-27 06:20:02'
$ws.Range("E45").Value = '2026-02-27 06:20:05'
$ws.Range("N45").Value = '5.0 °C 5:34 TU'
$ws.Range("O45").Value = '6.8 °C'
$ws.Range("E46").Value = '2026-02-27 06:20:07'
$ws.Range("J46").Value = '1026.1 hPa'
$ws.Range("N46").Value = '4.0 °C 5:54 TU'
$ws.Range("O46").Value = '6.5 °C'
